$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values (B2:E2)
$ws.Range("B2").Value = 10.713211367020678
$ws.Range("C2").Value = 8.8538691746208524
$ws.Range("D2").Value = 8.5158951133635679
$ws.Range("E2").Value = -0.18850371322243292

# Update row 3 values (B3:E3)
$ws.Range("B3").Value = 28.480680972077906
$ws.Range("C3").Value = 2.7472038679327975
$ws.Range("D3").Value = 1.3495636257163581
$ws.Range("E3").Value = 2.2223642054436512

# Update the selection to match the new selected range
$ws.Range("B1:E3").Select()
